$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Arkusz1")
$ws2 = $wb.Worksheets.Item("Arkusz2")

# Swap the "GNs" / "Gz" column headers in the three repeated blocks on Arkusz1 row 4
$ws1.Range("F4").Value = "Gz"
$ws1.Range("G4").Value = "GNs"
$ws1.Range("M4").Value = "Gz"
$ws1.Range("N4").Value = "GNs"
$ws1.Range("T4").Value = "Gz"
$ws1.Range("U4").Value = "GNs"

# Re-type "Razem" (capitalised) over the old lower-case "razem" labels
$ws1.Range("K4").Value = "Razem"
$ws1.Range("R4").Value = "Razem"
$ws1.Range("Y4").Value = "Razem"

# A stray keystroke landed in E9
$ws1.Range("E9").Value = "``"

# Fix up the judge-name header on Arkusz2: real line break instead of run of spaces
$ws2.Range("B3").Value = "Nazwisko i imię sędziego `nwg funkcji w wydziale"

# Leave the selections where the author left them, Arkusz2 ends up the active sheet
$ws1.Range("E9").Select()
$ws2.Range("B3:B4").Select()
